$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (31 and 32) following the existing pattern of the table,
# for two new Mac-Addresses (machine records).
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = 10030
$ws.Cells.Item(31, 4).Value = "eng"
$ws.Cells.Item(31, 5).Value = $true
$ws.Cells.Item(31, 6).Value = "superadmin"
$ws.Cells.Item(31, 7).Value = "now()"

$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = 10031
$ws.Cells.Item(32, 4).Value = "eng"
$ws.Cells.Item(32, 5).Value = $true
$ws.Cells.Item(32, 6).Value = "superadmin"
$ws.Cells.Item(32, 7).Value = "now()"

# Match the final selection/view state from the diff.
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 25
